$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column Q (year 2020), mirroring column P's formatting for each row
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)

# Fill in the 2020 values for the new column
$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 2
$ws.Range("Q6").Value = 0.3
$ws.Range("Q7").Value = 0.1
$ws.Range("Q8").Value = 4.3

$excel.CutCopyMode = 0

# Update the active selection, as captured in the workbook view
$ws.Range("O12").Select()
